# "Generate Report for Archive"
# - Update status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview sheet's zh-cn/de-de status columns, plus the per-locale
#   status column on the "zh-cn" and "de-de" sheets).
# - Shrink the now-narrower "Status" columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: column C (Status), rows 2-3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet: column C (Status), rows 2-3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Re-fit the columns that held the status text now that it is shorter.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
